$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.812.00"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "3.415.82"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'231.54"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "'619.74"
$ws.Range("E6").Value = "  -3.63%  "
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("D8").Value = "'0.391"
$ws.Range("E8").Value = "  -3.72%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("D11").Value = "3.413.37"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'42.93"
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("D13").Value = "'0.197"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Value = "'6.23"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "4.060.89"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "92.714.79"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "'8.09"
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("D19").Value = "3.415.33"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'17.88"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "'11.55"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'497.44"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'3.31"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").Value = "'0.438"
$ws.Range("E24").Value = "  -12.62%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'6.52"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000183"
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("D27").Value = "'91.24"
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("D28").Value = "'11.97"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "3.600.49"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'11.29"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'2.71"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("E33").Value = "  -3.47%  "
$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "'0.172"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("D36").Value = "'29.63"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "'0.540"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("D38").Value = "'553.67"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'7.47"
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.150"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'1.39"
$ws.Range("E42").Value = "  -4.66%  "
$ws.Range("D43").Value = "'0.912"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").Value = "'1.72"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").Value = "'3.71"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'23.65"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").Value = "'0.0406"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "'53.09"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").Value = "'2.10"
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("D51").Value = "'1.11"
$ws.Range("E51").Value = "  +16.14%  "
